# Leaderboard update: add "Inpromptu Rapid Fire" and "Chess" events, recompute totals,
# and add a bordered placeholder box (rank/house-logo slot) next to the new "Inpromptu
# Rapid Fire" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the old "total" row (row 11) completely so the sheet is contiguous
#    again once the new rows are inserted right after the existing events.
# ---------------------------------------------------------------------------
$ws.Range("A11:E11").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2. Row 7 - "Inpromptu Rapid Fire"
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Inpromptu Rapid`nFire"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 1

$ws.Range("A7:E7").Font.Name = "Arial"

# Placeholder box for rank + house logo (F7:I7), with G7:H7 merged.
$rF = $ws.Range("F7").Borders.Item(10)
$rF.Color = 0
$rF.LineStyle = 1

$gT = $ws.Range("G7").Borders.Item(8)
$gT.Color = 0
$gT.LineStyle = 1
$gB = $ws.Range("G7").Borders.Item(9)
$gB.Color = 0
$gB.LineStyle = 1

$hT = $ws.Range("H7").Borders.Item(8)
$hT.Color = 0
$hT.LineStyle = 1
$hB = $ws.Range("H7").Borders.Item(9)
$hB.Color = 0
$hB.LineStyle = 1
$hR = $ws.Range("H7").Borders.Item(10)
$hR.Color = 0
$hR.LineStyle = 1

$iT = $ws.Range("I7").Borders.Item(8)
$iT.Color = 0
$iT.LineStyle = 1
$iB = $ws.Range("I7").Borders.Item(9)
$iB.Color = 0
$iB.LineStyle = 1
$iR = $ws.Range("I7").Borders.Item(10)
$iR.Color = 0
$iR.LineStyle = 1

$ws.Range("G7:H7").Merge()

# Rest of the row keeps the same Arial styling out to column Z (mirrors the
# formatting Google Sheets paints across the full new row).
$ws.Range("F7:Z7").Font.Name = "Arial"

# ---------------------------------------------------------------------------
# 3. Row 8 - "Chess"
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Chess"
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1

$ws.Range("A8:E8").Font.Name = "Arial"

# ---------------------------------------------------------------------------
# 4. Row 9 - recomputed "total" row
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "total"
$ws.Range("B9").Value = 28
$ws.Range("C9").Value = 17
$ws.Range("D9").Value = 18
$ws.Range("E9").Value = 14

$ws.Range("B9:E9").HorizontalAlignment = -4152

Write-Host "leaderboard rows rebuilt: Inpromptu Rapid Fire, Chess, totals"
